# Strip the stray trailing "16" that was appended to every Bible reference
# in column A (e.g. "Ruth 1:116" -> "Ruth 1:1"), making the sheet human
# readable again. Header row (row 1) and column B are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2

    if ($null -ne $value -and $value.ToString().EndsWith("16")) {
        $newValue = $value.ToString().Substring(0, $value.ToString().Length - 2)
        $cell.Value2 = $newValue
    }
}
